$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 329.25
$ws.Range("I6").Value = 305.66666
$ws.Range("J6").Value = 400
$ws.Range("K6").Value = 916.9999799999999
$ws.Range("L6").Value = 1200
$ws.Range("M6").Value = -804.9999799999999
$ws.Range("N6").Value = -1424

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H10").Value = 3000
$ws.Range("I10").Value = 3000
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 3000
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = -2707

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 500
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 500
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 1500
$ws.Range("M31").ClearContents()
$ws.Range("N31").Value = -1960

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 2910.4482
$ws.Range("I74").Value = 2544.0833
$ws.Range("J74").Value = 3169.0588
$ws.Range("K74").Value = 2544.0833
$ws.Range("L74").Value = 3169.0588
$ws.Range("M74").Value = -1608.0833
$ws.Range("N74").Value = -5041.0588

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H77").Value = 2910.4482
$ws.Range("I77").Value = 2544.0833
$ws.Range("J77").Value = 3169.0588
$ws.Range("K77").Value = 12720.4165
$ws.Range("L77").Value = 15845.294
$ws.Range("M77").Value = -8040.416499999999
$ws.Range("N77").Value = -25205.294

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H123").Value = 70000
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 70000
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 70000
$ws.Range("N123").Value = -79800

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2599.2622
$ws.Range("I137").Value = 2437.5454
$ws.Range("J137").Value = 3017.8235
$ws.Range("K137").Value = 7312.6362
$ws.Range("L137").Value = 9053.470499999999
$ws.Range("M137").Value = -4762.6362
$ws.Range("N137").Value = -14153.4705

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 3199.625
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 3199.625
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 3199.625
$ws.Range("M74").ClearContents()
$ws.Range("N74").Value = -4947.625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 3199.625
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 3199.625
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 15998.125
$ws.Range("M77").ClearContents()
$ws.Range("N77").Value = -24734.125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1173.8334
$ws.Range("I122").Value = 905.2083
$ws.Range("J122").Value = 2248.3333
$ws.Range("K122").Value = 2715.6249
$ws.Range("L122").Value = 6744.999899999999
$ws.Range("M122").Value = -265.6248999999998
$ws.Range("N122").Value = -11644.9999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1106
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 1106
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 1106
$ws.Range("M31").ClearContents()
$ws.Range("N31").Value = -1696

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1106
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 1106
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 1106
$ws.Range("M34").ClearContents()
$ws.Range("N34").Value = -1510

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1583
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 1583
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 1583
$ws.Range("M58").ClearContents()
$ws.Range("N58").Value = -1989

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 1583
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 1583
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 4749
$ws.Range("M136").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 191.71428
$ws.Range("I6").Value = 140.33333
$ws.Range("J6").Value = 500
$ws.Range("K6").Value = 420.99999
$ws.Range("L6").Value = 1500
$ws.Range("M6").Value = -307.99999
$ws.Range("N6").Value = -1726

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 157172660
$ws.Range("I9").Value = 100000000
$ws.Range("J9").Value = 166701440
$ws.Range("K9").Value = 300000000
$ws.Range("L9").Value = 500104320
$ws.Range("M9").Value = -299999776
$ws.Range("N9").Value = -500104768

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 2000
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 2000
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 6000
$ws.Range("N64").Value = -6540

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H67").Value = 2000
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 2000
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 6000
$ws.Range("N67").Value = -7872

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 8595062
$ws.Range("I131").Value = 25003346
$ws.Range("J131").Value = 4333169.5
$ws.Range("K131").Value = 75010038
$ws.Range("L131").Value = 12999508.5
$ws.Range("M131").Value = -75004998
$ws.Range("N131").Value = -13009588.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1413.1923
$ws.Range("I102").Value = 1149.8823
$ws.Range("J102").Value = 1910.5555
$ws.Range("K102").Value = 1149.8823
$ws.Range("L102").Value = 1910.5555
$ws.Range("M102").Value = 472.1177
$ws.Range("N102").Value = -5154.5555

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H103").Value = 42300
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 42300
$ws.Range("K103").Value = 0
$ws.Range("L103").Value = 42300
$ws.Range("N103").Value = -44644

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 10001871
$ws.Range("I126").Value = 1185.3334
$ws.Range("J126").Value = 25002900
$ws.Range("K126").Value = 3556.0002
$ws.Range("L126").Value = 75008700
$ws.Range("M126").Value = -1086.0002
$ws.Range("N126").Value = -75013640

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2534415.5
$ws.Range("I22").Value = 2534415.5
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 2534415.5
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -2534120.5
$ws.Range("N22").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 2534415.5
$ws.Range("I27").Value = 2534415.5
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 2534415.5
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = -2534308.5
$ws.Range("N27").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 27779700
$ws.Range("I40").Value = 2161.5
$ws.Range("J40").Value = 250000000
$ws.Range("K40").Value = 2161.5
$ws.Range("L40").Value = 250000000
$ws.Range("M40").Value = -2025.5
$ws.Range("N40").Value = -250000272

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1151
$ws.Range("I68").Value = 1131.6154
$ws.Range("J68").Value = 1403
$ws.Range("K68").Value = 1131.6154
$ws.Range("L68").Value = 1403
$ws.Range("M68").Value = -382.6153999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H70").Value = 31813.334
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 31813.334
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 31813.334
$ws.Range("N70").Value = -32353.334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 1151
$ws.Range("I71").Value = 1131.6154
$ws.Range("J71").Value = 1403
$ws.Range("K71").Value = 5658.076999999999
$ws.Range("L71").Value = 7015
$ws.Range("M71").Value = -1914.076999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H73").Value = 31813.334
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 31813.334
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 31813.334
$ws.Range("N73").Value = -33685.334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 5861.524
$ws.Range("I122").Value = 8249.429
$ws.Range("J122").Value = 1085.7142
$ws.Range("K122").Value = 24748.287
$ws.Range("L122").Value = 3257.1426
$ws.Range("M122").Value = -22298.287
$ws.Range("N122").Value = -8157.142599999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 16906.941
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 16906.941
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 50720.823
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -55780.823

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4294.2
$ws.Range("I62").Value = 4104.6665
$ws.Range("J62").Value = 6000
$ws.Range("K62").Value = 4104.6665
$ws.Range("L62").Value = 6000
$ws.Range("M62").Value = -3480.6665

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 4294.2
$ws.Range("I65").Value = 4104.6665
$ws.Range("J65").Value = 6000
$ws.Range("K65").Value = 20523.3325
$ws.Range("L65").Value = 30000
$ws.Range("M65").Value = -17403.3325

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H82").Value = 59800
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 59800
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 59800
$ws.Range("N82").Value = -60566

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H85").Value = 59800
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 59800
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 59800
$ws.Range("N85").Value = -62452

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H92").Value = 51066.668
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 51066.668
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 51066.668
$ws.Range("N92").Value = -56058.668

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 27937.37
$ws.Range("I122").Value = 34340.668
$ws.Range("J122").Value = 3925
$ws.Range("K122").Value = 103022.004
$ws.Range("L122").Value = 11775
$ws.Range("M122").Value = -100572.004
$ws.Range("N122").Value = -16675

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2649.5557
$ws.Range("I126").Value = 2107.9092
$ws.Range("J126").Value = 3500.7144
$ws.Range("K126").Value = 6323.7276
$ws.Range("L126").Value = 10502.1432
$ws.Range("M126").Value = -3853.7276
$ws.Range("N126").Value = -15442.1432

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 13600.4
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 13600.4
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 40801.2
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -45861.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 14288312
$ws.Range("I136").Value = 35716224
$ws.Range("J136").Value = 3037.1428
$ws.Range("K136").Value = 107148672
$ws.Range("L136").Value = 9111.428400000001
$ws.Range("M136").Value = -107146122
